# Auto-generated edit script: applies the Sheets/Seraph_Profits.xlsx OOXML diff
# to the already-open workbook. 491 literal cell-value changes across 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); 3 of them remove a cell entirely
# (ClearContents), matching the source diff where an <c> element is deleted
# rather than its <v> merely changed.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (138 cell changes) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1101.5555 # H2 (was 1160.6)
$ws.Cells.Item(2, 9).Value = 1202.75 # I2 (was 2000)
$ws.Cells.Item(2, 10).Value = 1020.6 # J2 (was 950.75)
$ws.Cells.Item(2, 11).Value = 1202.75 # K2 (was 2000)
$ws.Cells.Item(2, 12).Value = 1020.6 # L2 (was 950.75)
$ws.Cells.Item(2, 13).Value = -1089.75 # M2 (was -1887)
$ws.Cells.Item(2, 14).Value = -1246.6 # N2 (was -1176.75)
$ws.Cells.Item(11, 8).Value = 438.3158 # H11 (was 421.25)
$ws.Cells.Item(11, 9).Value = 438.3158 # I11 (was 421.25)
$ws.Cells.Item(11, 11).Value = 438.3158 # K11 (was 421.25)
$ws.Cells.Item(11, 13).Value = -298.3158 # M11 (was -281.25)
$ws.Cells.Item(33, 8).Value = 443.4 # H33 (was 504.9)
$ws.Cells.Item(33, 9).Value = 131 # I33 (was 157.71428)
$ws.Cells.Item(33, 10).Value = 1068.2 # J33 (was 1315)
$ws.Cells.Item(33, 11).Value = 131 # K33 (was 157.71428)
$ws.Cells.Item(33, 12).Value = 1068.2 # L33 (was 1315)
$ws.Cells.Item(33, 13).Value = 98 # M33 (was 71.28572)
$ws.Cells.Item(33, 14).Value = -1526.2 # N33 (was -1773)
$ws.Cells.Item(38, 8).Value = 250.6 # H38 (was 294.5)
$ws.Cells.Item(38, 9).Value = 250.6 # I38 (was 294.5)
$ws.Cells.Item(38, 11).Value = 751.8 # K38 (was 883.5)
$ws.Cells.Item(38, 13).Value = -379.8 # M38 (was -511.5)
$ws.Cells.Item(41, 8).Value = 337.25 # H41 (was 362.33334)
$ws.Cells.Item(41, 9).Value = 322.5 # I41 (was 347.33334)
$ws.Cells.Item(41, 10).Value = 411 # J41 (was 407.33334)
$ws.Cells.Item(41, 11).Value = 322.5 # K41 (was 347.33334)
$ws.Cells.Item(41, 12).Value = 411 # L41 (was 407.33334)
$ws.Cells.Item(41, 13).Value = 117.5 # M41 (was 92.66665999999998)
$ws.Cells.Item(41, 14).Value = -1291 # N41 (was -1287.33334)
$ws.Cells.Item(43, 8).Value = 2000 # H43 (was 1997.5)
$ws.Cells.Item(43, 9).Value = 2000 # I43 (was 1997.5)
$ws.Cells.Item(43, 11).Value = 2000 # K43 (was 1997.5)
$ws.Cells.Item(43, 13).Value = -1931 # M43 (was -1928.5)
$ws.Cells.Item(48, 8).Value = 6500 # H48 (was 1000)
$ws.Cells.Item(48, 9).Value = 12000 # I48 (was 0)
$ws.Cells.Item(48, 11).Value = 36000 # K48 (was 0)
$ws.Cells.Item(48, 13).Value = -35708 # M48 (was empty)
$ws.Cells.Item(52, 8).Value = 428.1111 # H52 (was 331.75)
$ws.Cells.Item(52, 10).Value = 742.8 # J52 (was 628.75)
$ws.Cells.Item(52, 12).Value = 2228.4 # L52 (was 1886.25)
$ws.Cells.Item(52, 14).Value = -2548.4 # N52 (was -2206.25)
$ws.Cells.Item(56, 8).Value = 6500 # H56 (was 1000)
$ws.Cells.Item(56, 9).Value = 12000 # I56 (was 0)
$ws.Cells.Item(56, 11).Value = 36000 # K56 (was 0)
$ws.Cells.Item(56, 13).Value = -35466 # M56 (was empty)
$ws.Cells.Item(58, 8).Value = 2431.6 # H58 (was 2071.5)
$ws.Cells.Item(58, 9).Value = 62 # I58 (was 107.25)
$ws.Cells.Item(58, 10).Value = 4011.3333 # J58 (was 6000)
$ws.Cells.Item(58, 11).Value = 186 # K58 (was 321.75)
$ws.Cells.Item(58, 12).Value = 12033.9999 # L58 (was 18000)
$ws.Cells.Item(58, 13).Value = -36 # M58 (was -171.75)
$ws.Cells.Item(58, 14).Value = -12333.9999 # N58 (was -18300)
$ws.Cells.Item(76, 8).Value = 6218.3335 # H76 (was 5494.4165)
$ws.Cells.Item(76, 9).Value = 3001 # I76 (was 2500)
$ws.Cells.Item(76, 10).Value = 6620.5 # J76 (was 6093.3)
$ws.Cells.Item(76, 11).Value = 3001 # K76 (was 2500)
$ws.Cells.Item(76, 12).Value = 6620.5 # L76 (was 6093.3)
$ws.Cells.Item(76, 13).Value = -2686 # M76 (was -2185)
$ws.Cells.Item(76, 14).Value = -7250.5 # N76 (was -6723.3)
$ws.Cells.Item(79, 8).Value = 6218.3335 # H79 (was 5494.4165)
$ws.Cells.Item(79, 9).Value = 3001 # I79 (was 2500)
$ws.Cells.Item(79, 10).Value = 6620.5 # J79 (was 6093.3)
$ws.Cells.Item(79, 11).Value = 3001 # K79 (was 2500)
$ws.Cells.Item(79, 12).Value = 6620.5 # L79 (was 6093.3)
$ws.Cells.Item(79, 13).Value = -1909 # M79 (was -1408)
$ws.Cells.Item(79, 14).Value = -8804.5 # N79 (was -8277.299999999999)
$ws.Cells.Item(86, 8).Value = 2939.0625 # H86 (was 3289)
$ws.Cells.Item(86, 9).Value = 1468 # I86 (was 1859.4)
$ws.Cells.Item(86, 11).Value = 1468 # K86 (was 1859.4)
$ws.Cells.Item(86, 13).Value = -345 # M86 (was -736.4000000000001)
$ws.Cells.Item(89, 8).Value = 2939.0625 # H89 (was 3289)
$ws.Cells.Item(89, 9).Value = 1468 # I89 (was 1859.4)
$ws.Cells.Item(89, 11).Value = 7340 # K89 (was 9297)
$ws.Cells.Item(89, 13).Value = -1724 # M89 (was -3681)
$ws.Cells.Item(98, 8).Value = 1482.4783 # H98 (was 1682.85)
$ws.Cells.Item(98, 9).Value = 1424.619 # I98 (was 1637.6111)
$ws.Cells.Item(98, 11).Value = 1424.619 # K98 (was 1637.6111)
$ws.Cells.Item(98, 13).Value = 73.38100000000009 # M98 (was -139.6111000000001)
$ws.Cells.Item(100, 8).Value = 1115.6957 # H100 (was 1236.65)
$ws.Cells.Item(100, 9).Value = 763.05 # I100 (was 877.4666999999999)
$ws.Cells.Item(100, 10).Value = 3466.6667 # J100 (was 2314.2)
$ws.Cells.Item(100, 11).Value = 763.05 # K100 (was 877.4666999999999)
$ws.Cells.Item(100, 12).Value = 3466.6667 # L100 (was 2314.2)
$ws.Cells.Item(100, 13).Value = -222.05 # M100 (was -336.4666999999999)
$ws.Cells.Item(100, 14).Value = -4548.6667 # N100 (was -3396.2)
$ws.Cells.Item(101, 8).Value = 0 # H101 (was 900)
$ws.Cells.Item(101, 9).Value = 0 # I101 (was 900)
$ws.Cells.Item(101, 11).Value = 0 # K101 (was 2700)
$ws.Cells.Item(101, 13).ClearContents() # M101 (was -1078 -> removed)
$ws.Cells.Item(104, 8).Value = 125.75 # H104 (was 138)
$ws.Cells.Item(104, 10).Value = 89.5 # J104 (was 90)
$ws.Cells.Item(104, 12).Value = 268.5 # L104 (was 270)
$ws.Cells.Item(104, 14).Value = -3762.5 # N104 (was -3764)
$ws.Cells.Item(113, 9).Value = 5302.5 # I113 (was 5153)
$ws.Cells.Item(113, 10).Value = 7276.5 # J113 (was 9998)
$ws.Cells.Item(113, 11).Value = 5302.5 # K113 (was 5153)
$ws.Cells.Item(113, 12).Value = 7276.5 # L113 (was 9998)
$ws.Cells.Item(113, 13).Value = -2048.5 # M113 (was -1899)
$ws.Cells.Item(113, 14).Value = -13784.5 # N113 (was -16506)
$ws.Cells.Item(116, 8).Value = 5394.2856 # H116 (was 5660)
$ws.Cells.Item(116, 9).Value = 3952 # I116 (was 3990)
$ws.Cells.Item(116, 11).Value = 3952 # K116 (was 3990)
$ws.Cells.Item(116, 13).Value = -510 # M116 (was -548)
$ws.Cells.Item(118, 8).Value = 284.75 # H118 (was 309.75)
$ws.Cells.Item(118, 9).Value = 284.75 # I118 (was 309.75)
$ws.Cells.Item(118, 11).Value = 854.25 # K118 (was 929.25)
$ws.Cells.Item(118, 13).Value = 802.75 # M118 (was 727.75)
$ws.Cells.Item(122, 8).Value = 1482.4783 # H122 (was 1682.85)
$ws.Cells.Item(122, 9).Value = 1424.619 # I122 (was 1637.6111)
$ws.Cells.Item(122, 11).Value = 4273.857 # K122 (was 4912.8333)
$ws.Cells.Item(122, 13).Value = -1823.857 # M122 (was -2462.8333)
$ws.Cells.Item(127, 8).Value = 1214.1428 # H127 (was 1399.6666)
$ws.Cells.Item(127, 9).Value = 900 # I127 (was 1099.75)
$ws.Cells.Item(127, 11).Value = 2700 # K127 (was 3299.25)
$ws.Cells.Item(127, 13).Value = 2260 # M127 (was 1660.75)
$ws.Cells.Item(129, 8).Value = 548 # H129 (was 529.25)
$ws.Cells.Item(129, 9).Value = 548 # I129 (was 529.25)
$ws.Cells.Item(129, 11).Value = 1644 # K129 (was 1587.75)
$ws.Cells.Item(129, 13).Value = 3356 # M129 (was 3412.25)
$ws.Cells.Item(132, 8).Value = 3293.25 # H132 (was 3330)
$ws.Cells.Item(132, 9).Value = 3236.9443 # I132 (was 3274.7058)
$ws.Cells.Item(132, 11).Value = 9710.832900000001 # K132 (was 9824.117400000001)
$ws.Cells.Item(132, 13).Value = -7180.832900000001 # M132 (was -7294.117400000001)
$ws.Cells.Item(137, 8).Value = 11914 # H137 (was 9710.888999999999)
$ws.Cells.Item(137, 9).Value = 9733 # I137 (was 7799.75)
$ws.Cells.Item(137, 11).Value = 29199 # K137 (was 23399.25)
$ws.Cells.Item(137, 13).Value = -26649 # M137 (was -20849.25)
$ws.Cells.Item(138, 8).Value = 7512.57 # H138 (was 7739.108)
$ws.Cells.Item(138, 9).Value = 7314.0835 # I138 (was 7545.2173)
$ws.Cells.Item(138, 10).Value = 7589.4033 # J138 (was 7826.549)
$ws.Cells.Item(138, 11).Value = 21942.2505 # K138 (was 22635.6519)
$ws.Cells.Item(138, 12).Value = 22768.2099 # L138 (was 23479.647)
$ws.Cells.Item(138, 13).Value = -16802.2505 # M138 (was -17495.6519)
$ws.Cells.Item(138, 14).Value = -33048.2099 # N138 (was -33759.647)
$ws.Cells.Item(141, 8).Value = 735.5714 # H141 (was 800)
$ws.Cells.Item(141, 9).Value = 729.8 # I141 (was 825)
$ws.Cells.Item(141, 11).Value = 2189.4 # K141 (was 2475)
$ws.Cells.Item(141, 13).Value = 2990.6 # M141 (was 2705)

# ---- Sheet: ARM (35 cell changes) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 25708.607 # H32 (was 20305.945)
$ws.Cells.Item(32, 9).Value = 20787.7 # I32 (was 12597.685)
$ws.Cells.Item(32, 11).Value = 20787.7 # K32 (was 12597.685)
$ws.Cells.Item(32, 13).Value = -20500.7 # M32 (was -12310.685)
$ws.Cells.Item(61, 8).Value = 2337.3333 # H61 (was 2078.5715)
$ws.Cells.Item(61, 9).Value = 2337.3333 # I61 (was 2078.5715)
$ws.Cells.Item(61, 11).Value = 2337.3333 # K61 (was 2078.5715)
$ws.Cells.Item(61, 13).Value = -2125.3333 # M61 (was -1866.5715)
$ws.Cells.Item(97, 8).Value = 964.94116 # H97 (was 952.44446)
$ws.Cells.Item(97, 9).Value = 1010.3571 # I97 (was 992.3333)
$ws.Cells.Item(97, 11).Value = 1010.3571 # K97 (was 992.3333)
$ws.Cells.Item(97, 13).Value = -514.3570999999999 # M97 (was -496.3333)
$ws.Cells.Item(102, 8).Value = 1919.4783 # H102 (was 2002.2858)
$ws.Cells.Item(102, 9).Value = 1639 # I102 (was 1754.4445)
$ws.Cells.Item(102, 10).Value = 2135.2307 # J102 (was 2188.1667)
$ws.Cells.Item(102, 11).Value = 1639 # K102 (was 1754.4445)
$ws.Cells.Item(102, 12).Value = 2135.2307 # L102 (was 2188.1667)
$ws.Cells.Item(102, 13).Value = -17 # M102 (was -132.4445000000001)
$ws.Cells.Item(102, 14).Value = -5379.2307 # N102 (was -5432.1667)
$ws.Cells.Item(109, 8).Value = 0 # H109 (was 75188.5)
$ws.Cells.Item(109, 10).Value = 0 # J109 (was 75188.5)
$ws.Cells.Item(109, 12).Value = 0 # L109 (was 75188.5)
$ws.Cells.Item(109, 14).ClearContents() # N109 (was -77962.5 -> removed)
$ws.Cells.Item(122, 8).Value = 670933.8 # H122 (was 718571.9399999999)
$ws.Cells.Item(122, 10).Value = 4751.125 # J122 (was 4858.4287)
$ws.Cells.Item(122, 12).Value = 14253.375 # L122 (was 14575.2861)
$ws.Cells.Item(122, 14).Value = -19153.375 # N122 (was -19475.2861)
$ws.Cells.Item(132, 8).Value = 6203 # H132 (was 2847.2)
$ws.Cells.Item(132, 9).Value = 3345.5 # I132 (was 1744.7693)
$ws.Cells.Item(132, 11).Value = 10036.5 # K132 (was 5234.3079)
$ws.Cells.Item(132, 13).Value = -7506.5 # M132 (was -2704.3079)
$ws.Cells.Item(136, 8).Value = 2337.3333 # H136 (was 2078.5715)
$ws.Cells.Item(136, 9).Value = 2337.3333 # I136 (was 2078.5715)
$ws.Cells.Item(136, 11).Value = 7011.999899999999 # K136 (was 6235.7145)
$ws.Cells.Item(136, 13).Value = -4461.999899999999 # M136 (was -3685.7145)

# ---- Sheet: BSM (45 cell changes) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2771.2778 # H20 (was 2837.7222)
$ws.Cells.Item(20, 10).Value = 2870.4443 # J20 (was 3003.3333)
$ws.Cells.Item(20, 12).Value = 2870.4443 # L20 (was 3003.3333)
$ws.Cells.Item(20, 14).Value = -3364.4443 # N20 (was -3497.3333)
$ws.Cells.Item(22, 8).Value = 459.57895 # H22 (was 489.05264)
$ws.Cells.Item(22, 9).Value = 479.94446 # I22 (was 511.05554)
$ws.Cells.Item(22, 11).Value = 479.94446 # K22 (was 511.05554)
$ws.Cells.Item(22, 13).Value = -306.94446 # M22 (was -338.05554)
$ws.Cells.Item(86, 8).Value = 5386 # H86 (was 5454.3335)
$ws.Cells.Item(86, 9).Value = 4675.5 # I86 (was 4778)
$ws.Cells.Item(86, 11).Value = 4675.5 # K86 (was 4778)
$ws.Cells.Item(86, 13).Value = -3552.5 # M86 (was -3655)
$ws.Cells.Item(89, 8).Value = 5386 # H89 (was 5454.3335)
$ws.Cells.Item(89, 9).Value = 4675.5 # I89 (was 4778)
$ws.Cells.Item(89, 11).Value = 23377.5 # K89 (was 23890)
$ws.Cells.Item(89, 13).Value = -17761.5 # M89 (was -18274)
$ws.Cells.Item(94, 8).Value = 2422.5 # H94 (was 2507.6667)
$ws.Cells.Item(94, 9).Value = 2380.5 # I94 (was 2508)
$ws.Cells.Item(94, 10).Value = 2506.5 # J94 (was 2507)
$ws.Cells.Item(94, 11).Value = 2380.5 # K94 (was 2508)
$ws.Cells.Item(94, 12).Value = 2506.5 # L94 (was 2507)
$ws.Cells.Item(94, 13).Value = -1929.5 # M94 (was -2057)
$ws.Cells.Item(94, 14).Value = -3408.5 # N94 (was -3409)
$ws.Cells.Item(99, 8).Value = 1659.3 # H99 (was 1711)
$ws.Cells.Item(99, 9).Value = 1574.75 # I99 (was 1615.5)
$ws.Cells.Item(99, 11).Value = 1574.75 # K99 (was 1615.5)
$ws.Cells.Item(99, 13).Value = -76.75 # M99 (was -117.5)
$ws.Cells.Item(105, 8).Value = 4246.04 # H105 (was 4472.5)
$ws.Cells.Item(105, 9).Value = 3508.2942 # I105 (was 3615.0625)
$ws.Cells.Item(105, 10).Value = 5813.75 # J105 (was 6187.375)
$ws.Cells.Item(105, 11).Value = 3508.2942 # K105 (was 3615.0625)
$ws.Cells.Item(105, 12).Value = 5813.75 # L105 (was 6187.375)
$ws.Cells.Item(105, 13).Value = -1761.2942 # M105 (was -1868.0625)
$ws.Cells.Item(105, 14).Value = -9307.75 # N105 (was -9681.375)
$ws.Cells.Item(107, 8).Value = 6166 # H107 (was 6544.727)
$ws.Cells.Item(107, 9).Value = 4218.6665 # I107 (was 4496)
$ws.Cells.Item(107, 11).Value = 4218.6665 # K107 (was 4496)
$ws.Cells.Item(107, 13).Value = -2298.6665 # M107 (was -2576)
$ws.Cells.Item(134, 8).Value = 2882.2222 # H134 (was 2988.6667)
$ws.Cells.Item(134, 9).Value = 1433.75 # I134 (was 1496.1428)
$ws.Cells.Item(134, 10).Value = 4989.091 # J134 (was 5078.2)
$ws.Cells.Item(134, 11).Value = 4301.25 # K134 (was 4488.428400000001)
$ws.Cells.Item(134, 12).Value = 14967.273 # L134 (was 15234.6)
$ws.Cells.Item(134, 13).Value = -1766.25 # M134 (was -1953.428400000001)
$ws.Cells.Item(134, 14).Value = -20037.273 # N134 (was -20304.6)

# ---- Sheet: CRP (71 cell changes) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 724.5 # H22 (was 999.5)
$ws.Cells.Item(22, 9).Value = 724.5 # I22 (was 0)
$ws.Cells.Item(22, 10).Value = 0 # J22 (was 999.5)
$ws.Cells.Item(22, 11).Value = 724.5 # K22 (was 0)
$ws.Cells.Item(22, 12).Value = 0 # L22 (was 999.5)
$ws.Cells.Item(22, 13).Value = -374.5 # M22 (was empty)
$ws.Cells.Item(22, 14).ClearContents() # N22 (was -1699.5 -> removed)
$ws.Cells.Item(31, 8).Value = 4498.0347 # H31 (was 4549.4287)
$ws.Cells.Item(31, 9).Value = 3038.8823 # I31 (was 3240.7334)
$ws.Cells.Item(31, 10).Value = 6565.1665 # J31 (was 6059.4614)
$ws.Cells.Item(31, 11).Value = 3038.8823 # K31 (was 3240.7334)
$ws.Cells.Item(31, 12).Value = 6565.1665 # L31 (was 6059.4614)
$ws.Cells.Item(31, 13).Value = -2743.8823 # M31 (was -2945.7334)
$ws.Cells.Item(31, 14).Value = -7155.1665 # N31 (was -6649.4614)
$ws.Cells.Item(34, 8).Value = 4498.0347 # H34 (was 4549.4287)
$ws.Cells.Item(34, 9).Value = 3038.8823 # I34 (was 3240.7334)
$ws.Cells.Item(34, 10).Value = 6565.1665 # J34 (was 6059.4614)
$ws.Cells.Item(34, 11).Value = 3038.8823 # K34 (was 3240.7334)
$ws.Cells.Item(34, 12).Value = 6565.1665 # L34 (was 6059.4614)
$ws.Cells.Item(34, 13).Value = -2836.8823 # M34 (was -3038.7334)
$ws.Cells.Item(34, 14).Value = -6969.1665 # N34 (was -6463.4614)
$ws.Cells.Item(58, 8).Value = 4829.909 # H58 (was 4618.9165)
$ws.Cells.Item(58, 9).Value = 1448.8334 # I58 (was 1570.1428)
$ws.Cells.Item(58, 11).Value = 1448.8334 # K58 (was 1570.1428)
$ws.Cells.Item(58, 13).Value = -1245.8334 # M58 (was -1367.1428)
$ws.Cells.Item(86, 8).Value = 6696.522 # H86 (was 6654.9165)
$ws.Cells.Item(86, 9).Value = 3001.8 # I86 (was 3170.3125)
$ws.Cells.Item(86, 11).Value = 3001.8 # K86 (was 3170.3125)
$ws.Cells.Item(86, 13).Value = -1878.8 # M86 (was -2047.3125)
$ws.Cells.Item(89, 8).Value = 6696.522 # H89 (was 6654.9165)
$ws.Cells.Item(89, 9).Value = 3001.8 # I89 (was 3170.3125)
$ws.Cells.Item(89, 11).Value = 15009 # K89 (was 15851.5625)
$ws.Cells.Item(89, 13).Value = -9393 # M89 (was -10235.5625)
$ws.Cells.Item(99, 8).Value = 17153.23 # H99 (was 16643.143)
$ws.Cells.Item(99, 10).Value = 18873.375 # J99 (was 17888.777)
$ws.Cells.Item(99, 12).Value = 18873.375 # L99 (was 17888.777)
$ws.Cells.Item(99, 14).Value = -21869.375 # N99 (was -20884.777)
$ws.Cells.Item(105, 8).Value = 3274.5 # H105 (was 4372.5454)
$ws.Cells.Item(105, 9).Value = 1896.1 # I105 (was 2823.8572)
$ws.Cells.Item(105, 10).Value = 5571.8335 # J105 (was 7082.75)
$ws.Cells.Item(105, 11).Value = 1896.1 # K105 (was 2823.8572)
$ws.Cells.Item(105, 12).Value = 5571.8335 # L105 (was 7082.75)
$ws.Cells.Item(105, 13).Value = -149.0999999999999 # M105 (was -1076.8572)
$ws.Cells.Item(105, 14).Value = -9065.833500000001 # N105 (was -10576.75)
$ws.Cells.Item(107, 8).Value = 398.5 # H107 (was 420.55554)
$ws.Cells.Item(107, 9).Value = 314 # I107 (was 352)
$ws.Cells.Item(107, 11).Value = 314 # K107 (was 352)
$ws.Cells.Item(107, 13).Value = 1606 # M107 (was 1568)
$ws.Cells.Item(109, 8).Value = 64265 # H109 (was 64053.57)
$ws.Cells.Item(109, 10).Value = 64265 # J109 (was 64053.57)
$ws.Cells.Item(109, 12).Value = 64265 # L109 (was 64053.57)
$ws.Cells.Item(109, 14).Value = -66345 # N109 (was -66133.57000000001)
$ws.Cells.Item(126, 8).Value = 17153.23 # H126 (was 16643.143)
$ws.Cells.Item(126, 10).Value = 18873.375 # J126 (was 17888.777)
$ws.Cells.Item(126, 12).Value = 56620.125 # L126 (was 53666.33099999999)
$ws.Cells.Item(126, 14).Value = -61560.125 # N126 (was -58606.33099999999)
$ws.Cells.Item(132, 8).Value = 3412.6316 # H132 (was 2203.9412)
$ws.Cells.Item(132, 9).Value = 3622.3076 # I132 (was 2042.2858)
$ws.Cells.Item(132, 11).Value = 10866.9228 # K132 (was 6126.857400000001)
$ws.Cells.Item(132, 13).Value = -8336.9228 # M132 (was -3596.857400000001)
$ws.Cells.Item(134, 8).Value = 3801.75 # H134 (was 2336.973)
$ws.Cells.Item(134, 9).Value = 3043.375 # I134 (was 1746.2307)
$ws.Cells.Item(134, 10).Value = 4560.125 # J134 (was 3733.2727)
$ws.Cells.Item(134, 11).Value = 9130.125 # K134 (was 5238.6921)
$ws.Cells.Item(134, 12).Value = 13680.375 # L134 (was 11199.8181)
$ws.Cells.Item(134, 13).Value = -6595.125 # M134 (was -2703.6921)
$ws.Cells.Item(134, 14).Value = -18750.375 # N134 (was -16269.8181)
$ws.Cells.Item(136, 8).Value = 4829.909 # H136 (was 4618.9165)
$ws.Cells.Item(136, 9).Value = 1448.8334 # I136 (was 1570.1428)
$ws.Cells.Item(136, 11).Value = 4346.5002 # K136 (was 4710.428400000001)
$ws.Cells.Item(136, 13).Value = -1796.5002 # M136 (was -2160.428400000001)

# ---- Sheet: CUL (65 cell changes) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 233.66667 # H23 (was 258.91666)
$ws.Cells.Item(23, 9).Value = 293.33334 # I23 (was 180)
$ws.Cells.Item(23, 10).Value = 221.73334 # J23 (was 266.0909)
$ws.Cells.Item(23, 11).Value = 880.0000200000001 # K23 (was 540)
$ws.Cells.Item(23, 12).Value = 665.20002 # L23 (was 798.2727)
$ws.Cells.Item(23, 13).Value = -645.0000200000001 # M23 (was -305)
$ws.Cells.Item(23, 14).Value = -1135.20002 # N23 (was -1268.2727)
$ws.Cells.Item(34, 8).Value = 2466.3333 # H34 (was 3340)
$ws.Cells.Item(34, 9).Value = 1199.5 # I34 (was 2233.3333)
$ws.Cells.Item(34, 11).Value = 3598.5 # K34 (was 6699.999899999999)
$ws.Cells.Item(34, 13).Value = -3514.5 # M34 (was -6615.999899999999)
$ws.Cells.Item(39, 8).Value = 5000 # H39 (was 1666.6666)
$ws.Cells.Item(39, 10).Value = 15000 # J39 (was 0)
$ws.Cells.Item(39, 12).Value = 45000 # L39 (was 0)
$ws.Cells.Item(39, 14).Value = -45588 # N39 (was empty)
$ws.Cells.Item(55, 8).Value = 84899.5 # H55 (was 92872.17999999999)
$ws.Cells.Item(55, 9).Value = 250698.5 # I55 (was 200758.8)
$ws.Cells.Item(55, 10).Value = 2000 # J55 (was 2966.6667)
$ws.Cells.Item(55, 11).Value = 752095.5 # K55 (was 602276.3999999999)
$ws.Cells.Item(55, 12).Value = 6000 # L55 (was 8900.000100000001)
$ws.Cells.Item(55, 13).Value = -751918.5 # M55 (was -602099.3999999999)
$ws.Cells.Item(55, 14).Value = -6354 # N55 (was -9254.000100000001)
$ws.Cells.Item(107, 8).Value = 1445.5 # H107 (was 1692.3334)
$ws.Cells.Item(107, 9).Value = 2344.6 # I107 (was 2899.5)
$ws.Cells.Item(107, 10).Value = 946 # J107 (was 1088.75)
$ws.Cells.Item(107, 11).Value = 7033.799999999999 # K107 (was 8698.5)
$ws.Cells.Item(107, 12).Value = 2838 # L107 (was 3266.25)
$ws.Cells.Item(107, 13).Value = -5113.799999999999 # M107 (was -6778.5)
$ws.Cells.Item(107, 14).Value = -6678 # N107 (was -7106.25)
$ws.Cells.Item(113, 8).Value = 5084.467 # H113 (was 4415.048)
$ws.Cells.Item(113, 9).Value = 875 # I113 (was 662.5)
$ws.Cells.Item(113, 10).Value = 5385.143 # J113 (was 4810.0527)
$ws.Cells.Item(113, 11).Value = 2625 # K113 (was 1987.5)
$ws.Cells.Item(113, 12).Value = 16155.429 # L113 (was 14430.1581)
$ws.Cells.Item(113, 13).Value = -455 # M113 (was 182.5)
$ws.Cells.Item(113, 14).Value = -20495.429 # N113 (was -18770.1581)
$ws.Cells.Item(121, 8).Value = 2680.5557 # H121 (was 1938.125)
$ws.Cells.Item(121, 9).Value = 905.4 # I121 (was 917.6667)
$ws.Cells.Item(121, 10).Value = 4899.5 # J121 (was 4999.5)
$ws.Cells.Item(121, 11).Value = 2716.2 # K121 (was 2753.0001)
$ws.Cells.Item(121, 12).Value = 14698.5 # L121 (was 14998.5)
$ws.Cells.Item(121, 13).Value = -1406.2 # M121 (was -1443.0001)
$ws.Cells.Item(121, 14).Value = -17318.5 # N121 (was -17618.5)
$ws.Cells.Item(122, 8).Value = 2679.8 # H122 (was 3199.6)
$ws.Cells.Item(122, 9).Value = 2500 # I122 (was 2666.3333)
$ws.Cells.Item(122, 10).Value = 2799.6667 # J122 (was 3999.5)
$ws.Cells.Item(122, 11).Value = 22500 # K122 (was 23996.9997)
$ws.Cells.Item(122, 12).Value = 25197.0003 # L122 (was 35995.5)
$ws.Cells.Item(122, 13).Value = -20050 # M122 (was -21546.9997)
$ws.Cells.Item(122, 14).Value = -30097.0003 # N122 (was -40895.5)
$ws.Cells.Item(137, 8).Value = 7709.143 # H137 (was 7975.857)
$ws.Cells.Item(137, 9).Value = 7677.3335 # I137 (was 9000)
$ws.Cells.Item(137, 10).Value = 7733 # J137 (was 7566.2)
$ws.Cells.Item(137, 11).Value = 23032.0005 # K137 (was 27000)
$ws.Cells.Item(137, 12).Value = 23199 # L137 (was 22698.6)
$ws.Cells.Item(137, 13).Value = -17932.0005 # M137 (was -21900)
$ws.Cells.Item(137, 14).Value = -33399 # N137 (was -32898.6)
$ws.Cells.Item(139, 8).Value = 8643.909 # H139 (was 9259.5)
$ws.Cells.Item(139, 9).Value = 7155.2856 # I139 (was 7933.1665)
$ws.Cells.Item(139, 11).Value = 21465.8568 # K139 (was 23799.4995)
$ws.Cells.Item(139, 13).Value = -16325.8568 # M139 (was -18659.4995)
$ws.Cells.Item(140, 8).Value = 4798.1113 # H140 (was 5023)
$ws.Cells.Item(140, 9).Value = 4169 # I140 (was 4364)
$ws.Cells.Item(140, 11).Value = 12507 # K140 (was 13092)
$ws.Cells.Item(140, 13).Value = -7327 # M140 (was -7912)

# ---- Sheet: GSM (33 cell changes) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 192.14285 # H2 (was 172.34782)
$ws.Cells.Item(2, 9).Value = 32 # I2 (was 22.076923)
$ws.Cells.Item(2, 10).Value = 312.25 # J2 (was 367.7)
$ws.Cells.Item(2, 11).Value = 32 # K2 (was 22.076923)
$ws.Cells.Item(2, 12).Value = 312.25 # L2 (was 367.7)
$ws.Cells.Item(2, 13).Value = 81 # M2 (was 90.92307700000001)
$ws.Cells.Item(2, 14).Value = -538.25 # N2 (was -593.7)
$ws.Cells.Item(11, 8).Value = 3002493 # H11 (was 3275247.2)
$ws.Cells.Item(11, 10).Value = 2768.3333 # J11 (was 2840)
$ws.Cells.Item(11, 12).Value = 2768.3333 # L11 (was 2840)
$ws.Cells.Item(11, 14).Value = -3046.3333 # N11 (was -3118)
$ws.Cells.Item(97, 8).Value = 1906.091 # H97 (was 1969.7273)
$ws.Cells.Item(97, 9).Value = 1938.6842 # I97 (was 2012.3684)
$ws.Cells.Item(97, 11).Value = 1938.6842 # K97 (was 2012.3684)
$ws.Cells.Item(97, 13).Value = -1442.6842 # M97 (was -1516.3684)
$ws.Cells.Item(102, 8).Value = 4028.111 # H102 (was 4184.5)
$ws.Cells.Item(102, 9).Value = 2813.25 # I102 (was 2825.3333)
$ws.Cells.Item(102, 11).Value = 2813.25 # K102 (was 2825.3333)
$ws.Cells.Item(102, 13).Value = -1191.25 # M102 (was -1203.3333)
$ws.Cells.Item(126, 8).Value = 4835.3335 # H126 (was 4602.875)
$ws.Cells.Item(126, 9).Value = 4012 # I126 (was 3412)
$ws.Cells.Item(126, 10).Value = 5000 # J126 (was 4999.8335)
$ws.Cells.Item(126, 11).Value = 12036 # K126 (was 10236)
$ws.Cells.Item(126, 12).Value = 15000 # L126 (was 14999.5005)
$ws.Cells.Item(126, 13).Value = -9566 # M126 (was -7766)
$ws.Cells.Item(126, 14).Value = -19940 # N126 (was -19939.5005)
$ws.Cells.Item(132, 8).Value = 4667 # H132 (was 4600.6665)
$ws.Cells.Item(132, 9).Value = 3212.818 # I132 (was 3375.7144)
$ws.Cells.Item(132, 10).Value = 8666 # J132 (was 8888)
$ws.Cells.Item(132, 11).Value = 9638.454000000002 # K132 (was 10127.1432)
$ws.Cells.Item(132, 12).Value = 25998 # L132 (was 26664)
$ws.Cells.Item(132, 13).Value = -7108.454000000002 # M132 (was -7597.143199999999)
$ws.Cells.Item(132, 14).Value = -31058 # N132 (was -31724)

# ---- Sheet: LTW (52 cell changes) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 7000 # H13 (was 4000)
$ws.Cells.Item(13, 10).Value = 10000 # J13 (was 0)
$ws.Cells.Item(13, 12).Value = 10000 # L13 (was 0)
$ws.Cells.Item(13, 14).Value = -10280 # N13 (was empty)
$ws.Cells.Item(16, 8).Value = 2554.65 # H16 (was 3287.9333)
$ws.Cells.Item(16, 9).Value = 2887.625 # I16 (was 3524.5386)
$ws.Cells.Item(16, 10).Value = 1222.75 # J16 (was 1750)
$ws.Cells.Item(16, 11).Value = 2887.625 # K16 (was 3524.5386)
$ws.Cells.Item(16, 12).Value = 1222.75 # L16 (was 1750)
$ws.Cells.Item(16, 13).Value = -2717.625 # M16 (was -3354.5386)
$ws.Cells.Item(16, 14).Value = -1562.75 # N16 (was -2090)
$ws.Cells.Item(22, 8).Value = 4000 # H22 (was 300)
$ws.Cells.Item(22, 9).Value = 3000 # I22 (was 300)
$ws.Cells.Item(22, 10).Value = 5000 # J22 (was 0)
$ws.Cells.Item(22, 11).Value = 3000 # K22 (was 300)
$ws.Cells.Item(22, 12).Value = 5000 # L22 (was 0)
$ws.Cells.Item(22, 13).Value = -2705 # M22 (was -5)
$ws.Cells.Item(22, 14).Value = -5590 # N22 (was empty)
$ws.Cells.Item(27, 8).Value = 4000 # H27 (was 300)
$ws.Cells.Item(27, 9).Value = 3000 # I27 (was 300)
$ws.Cells.Item(27, 10).Value = 5000 # J27 (was 0)
$ws.Cells.Item(27, 11).Value = 3000 # K27 (was 300)
$ws.Cells.Item(27, 12).Value = 5000 # L27 (was 0)
$ws.Cells.Item(27, 13).Value = -2893 # M27 (was -193)
$ws.Cells.Item(27, 14).Value = -5214 # N27 (was empty)
$ws.Cells.Item(61, 8).Value = 5671.423 # H61 (was 5999.4585)
$ws.Cells.Item(61, 9).Value = 5628.609 # I61 (was 5999.4287)
$ws.Cells.Item(61, 11).Value = 5628.609 # K61 (was 5999.4287)
$ws.Cells.Item(61, 13).Value = -5426.609 # M61 (was -5797.4287)
$ws.Cells.Item(100, 8).Value = 4000.4443 # H100 (was 4643.2856)
$ws.Cells.Item(100, 9).Value = 3857.7144 # I100 (was 4583.8335)
$ws.Cells.Item(100, 10).Value = 4500 # J100 (was 5000)
$ws.Cells.Item(100, 11).Value = 3857.7144 # K100 (was 4583.8335)
$ws.Cells.Item(100, 12).Value = 4500 # L100 (was 5000)
$ws.Cells.Item(100, 13).Value = -3316.7144 # M100 (was -4042.8335)
$ws.Cells.Item(100, 14).Value = -5582 # N100 (was -6082)
$ws.Cells.Item(113, 8).Value = 5671.423 # H113 (was 5999.4585)
$ws.Cells.Item(113, 9).Value = 5628.609 # I113 (was 5999.4287)
$ws.Cells.Item(113, 11).Value = 5628.609 # K113 (was 5999.4287)
$ws.Cells.Item(113, 13).Value = -3458.609 # M113 (was -3829.4287)
$ws.Cells.Item(122, 8).Value = 9333.333000000001 # H122 (was 7099.8184)
$ws.Cells.Item(122, 9).Value = 8000 # I122 (was 5442.5713)
$ws.Cells.Item(122, 11).Value = 24000 # K122 (was 16327.7139)
$ws.Cells.Item(122, 13).Value = -21550 # M122 (was -13877.7139)
$ws.Cells.Item(132, 8).Value = 4756.3335 # H132 (was 4540.2666)
$ws.Cells.Item(132, 10).Value = 6319.091 # J132 (was 6500)
$ws.Cells.Item(132, 12).Value = 18957.273 # L132 (was 19500)
$ws.Cells.Item(132, 14).Value = -24017.273 # N132 (was -24560)
$ws.Cells.Item(136, 8).Value = 3184.1667 # H136 (was 3404)
$ws.Cells.Item(136, 9).Value = 3184.1667 # I136 (was 3404)
$ws.Cells.Item(136, 11).Value = 9552.500100000001 # K136 (was 10212)
$ws.Cells.Item(136, 13).Value = -7002.500100000001 # M136 (was -7662)

# ---- Sheet: WVR (52 cell changes) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(26, 8).Value = 529899.8 # H26 (was 653625)
$ws.Cells.Item(26, 9).Value = 37374.75 # I26 (was 38166.668)
$ws.Cells.Item(26, 11).Value = 37374.75 # K26 (was 38166.668)
$ws.Cells.Item(26, 13).Value = -37081.75 # M26 (was -37873.668)
$ws.Cells.Item(70, 8).Value = 58119 # H70 (was 57119)
$ws.Cells.Item(70, 9).Value = 55595 # I70 (was 52797.5)
$ws.Cells.Item(70, 10).Value = 58750 # J70 (was 60000)
$ws.Cells.Item(70, 11).Value = 55595 # K70 (was 52797.5)
$ws.Cells.Item(70, 12).Value = 58750 # L70 (was 60000)
$ws.Cells.Item(70, 13).Value = -55280 # M70 (was -52482.5)
$ws.Cells.Item(70, 14).Value = -59380 # N70 (was -60630)
$ws.Cells.Item(73, 8).Value = 58119 # H73 (was 57119)
$ws.Cells.Item(73, 9).Value = 55595 # I73 (was 52797.5)
$ws.Cells.Item(73, 10).Value = 58750 # J73 (was 60000)
$ws.Cells.Item(73, 11).Value = 55595 # K73 (was 52797.5)
$ws.Cells.Item(73, 12).Value = 58750 # L73 (was 60000)
$ws.Cells.Item(73, 13).Value = -54503 # M73 (was -51705.5)
$ws.Cells.Item(73, 14).Value = -60934 # N73 (was -62184)
$ws.Cells.Item(96, 8).Value = 2066.6667 # H96 (was 2500)
$ws.Cells.Item(96, 9).Value = 2066.6667 # I96 (was 2500)
$ws.Cells.Item(96, 11).Value = 2066.6667 # K96 (was 2500)
$ws.Cells.Item(96, 13).Value = -693.6667000000002 # M96 (was -1127)
$ws.Cells.Item(100, 8).Value = 2506.4285 # H100 (was 2595.8572)
$ws.Cells.Item(100, 9).Value = 3199 # I100 (was 3338.111)
$ws.Cells.Item(100, 11).Value = 6398 # K100 (was 6676.222)
$ws.Cells.Item(100, 13).Value = -5857 # M100 (was -6135.222)
$ws.Cells.Item(107, 8).Value = 2327 # H107 (was 1882.7142)
$ws.Cells.Item(107, 9).Value = 685 # I107 (was 743)
$ws.Cells.Item(107, 11).Value = 2055 # K107 (was 2229)
$ws.Cells.Item(107, 13).Value = -135 # M107 (was -309)
$ws.Cells.Item(122, 8).Value = 1525.4286 # H122 (was 1713)
$ws.Cells.Item(122, 9).Value = 1525.4286 # I122 (was 1713)
$ws.Cells.Item(122, 11).Value = 4576.2858 # K122 (was 5139)
$ws.Cells.Item(122, 13).Value = -2126.2858 # M122 (was -2689)
$ws.Cells.Item(126, 8).Value = 94898.91 # H126 (was 104369.3)
$ws.Cells.Item(126, 9).Value = 102888.8 # I126 (was 114299.22)
$ws.Cells.Item(126, 11).Value = 308666.4 # K126 (was 342897.66)
$ws.Cells.Item(126, 13).Value = -306196.4 # M126 (was -340427.66)
$ws.Cells.Item(132, 8).Value = 3024.5 # H132 (was 1440.909)
$ws.Cells.Item(132, 9).Value = 1350.5 # I132 (was 613.6667)
$ws.Cells.Item(132, 10).Value = 3582.5 # J132 (was 3213.5715)
$ws.Cells.Item(132, 11).Value = 4051.5 # K132 (was 1841.0001)
$ws.Cells.Item(132, 12).Value = 10747.5 # L132 (was 9640.7145)
$ws.Cells.Item(132, 13).Value = -1521.5 # M132 (was 688.9999)
$ws.Cells.Item(132, 14).Value = -15807.5 # N132 (was -14700.7145)
$ws.Cells.Item(136, 8).Value = 52255.9 # H136 (was 60950.117)
$ws.Cells.Item(136, 9).Value = 1153.75 # I136 (was 1281.7693)
$ws.Cells.Item(136, 10).Value = 256664.5 # J136 (was 254872.25)
$ws.Cells.Item(136, 11).Value = 3461.25 # K136 (was 3845.3079)
$ws.Cells.Item(136, 12).Value = 769993.5 # L136 (was 764616.75)
$ws.Cells.Item(136, 13).Value = -911.25 # M136 (was -1295.3079)
$ws.Cells.Item(136, 14).Value = -775093.5 # N136 (was -769716.75)
